$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------
# Apply cell formatting for the two new daily blocks by copying the
# formats from existing rows that already carry the right style mix
# (row 345 = date/Domm row, 346 = Meeting row w/ no A style,
#  347 = mid row w/ A:s2, 350 = mid row w/ A:s3/B:s3 incl totals row).
# -----------------------------------------------------------------

# Pattern "date row" (A:s2 B:s3 C:s3 D:s4) -> rows 354, 362
$ws.Range("A345:D345").Copy()
$ws.Range("A354:D354").PasteSpecial(-4122)
$ws.Range("A345:D345").Copy()
$ws.Range("A362:D362").PasteSpecial(-4122)

# Pattern "meeting row" (no A style, B:(none) C:s3 D:s4) -> rows 355, 363
$ws.Range("A346:D346").Copy()
$ws.Range("A355:D355").PasteSpecial(-4122)
$ws.Range("A346:D346").Copy()
$ws.Range("A363:D363").PasteSpecial(-4122)

# Pattern "A:s2, no B" -> rows 356, 364
$ws.Range("A347:D347").Copy()
$ws.Range("A356:D356").PasteSpecial(-4122)
$ws.Range("A347:D347").Copy()
$ws.Range("A364:D364").PasteSpecial(-4122)

# Pattern "A:s3 B:s3 C:s3 D:s4" (study / sub rows / total row) -> rows
# 357-360 and 365-368
$ws.Range("A350:D350").Copy()
$ws.Range("A357:D357").PasteSpecial(-4122)
$ws.Range("A350:D350").Copy()
$ws.Range("A358:D358").PasteSpecial(-4122)
$ws.Range("A350:D350").Copy()
$ws.Range("A359:D359").PasteSpecial(-4122)
$ws.Range("A350:D350").Copy()
$ws.Range("A360:D360").PasteSpecial(-4122)
$ws.Range("A350:D350").Copy()
$ws.Range("A365:D365").PasteSpecial(-4122)
$ws.Range("A350:D350").Copy()
$ws.Range("A366:D366").PasteSpecial(-4122)
$ws.Range("A350:D350").Copy()
$ws.Range("A367:D367").PasteSpecial(-4122)
$ws.Range("A350:D350").Copy()
$ws.Range("A368:D368").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# -----------------------------------------------------------------
# Block 1: Monday 2025-03-13 (serial 45729), rows 354-360
# -----------------------------------------------------------------
$ws.Range("A354").Value = 45729
$ws.Range("B354").Value = "Domm"
$ws.Range("D354").Value = 0.25

$ws.Range("B355").Value = "Meeting"
$ws.Range("C355").Value = "General Discussion"
$ws.Range("D355").Value = 0.25

$ws.Range("C356").Value = "Operational Training Reconsile"
$ws.Range("D356").Value = 0.5

$ws.Range("B357").Value = "Study"
$ws.Range("C357").Value = "Operational Training Practice "
$ws.Range("D357").Value = 2

$ws.Range("C358").Value = "Redis"
$ws.Range("D358").Value = 3

$ws.Range("C359").Value = "Documentation"
$ws.Range("D359").Value = 2

$ws.Range("B360").Value = "Total"
$ws.Range("D360").Formula = "=SUM(D353:D359)"

# -----------------------------------------------------------------
# Block 2: Friday 2025-03-17 (serial 45733), rows 362-368
# -----------------------------------------------------------------
$ws.Range("A362").Value = 45733
$ws.Range("B362").Value = "Domm"
$ws.Range("D362").Value = 0.25

$ws.Range("B363").Value = "Meeting"
$ws.Range("C363").Value = "General Discussion"
$ws.Range("D363").Value = 0.25

$ws.Range("C364").Value = "CM API Training + Discussion "
$ws.Range("D364").Value = 1

$ws.Range("B365").Value = "Study"
$ws.Range("C365").Value = "DevExtreme Documentation"
$ws.Range("D365").Value = 3

$ws.Range("C366").Value = "DevExtreme Project changes"
$ws.Range("D366").Value = 1

$ws.Range("C367").Value = "Read React Documentation"
$ws.Range("D367").Value = 2.5

$ws.Range("B368").Value = "Total"
$ws.Range("D368").Formula = "=SUM(D361:D367)"

# -----------------------------------------------------------------
# Update the view state to reflect the newly-added rows (best effort;
# mirrors the scroll/selection the author ended up with).
# -----------------------------------------------------------------
$ws.Range("A362:D369").Select()
